$wb = $excel.ActiveWorkbook

# Insert the new sheet right after "FTNC_Average_Demand151" (i.e. at the end)
$template = $wb.Worksheets.Item("FTNC_Average_Demand151")
$newSheet = $wb.Worksheets.Add($null, $template)
$newSheet.Name = "FTNC_Average_Demand152"

# Replicate the template's cell formatting (bold header / border style) onto the new sheet
$template.Range("A1:F2").Copy()
$newSheet.Range("A1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Header row
$newSheet.Range("B1").Value = "In-vehicle"
$newSheet.Range("C1").Value = "At-stop"
$newSheet.Range("D1").Value = "Extra"
$newSheet.Range("E1").Value = "Tardiness"
$newSheet.Range("F1").Value = "Total"

# Data row
$newSheet.Range("A2").Value = "FTNC_Average_Demand_15"
$newSheet.Range("B2").Value = 2444.36605568987
$newSheet.Range("C2").Value = 12957.52356681243
$newSheet.Range("D2").Value = 704.8661188217991
$newSheet.Range("E2").Value = 18.85729127237325
$newSheet.Range("F2").Value = 16125.61303259648
